# Auto-generated edit script: updates Leve profit/price figures across all 8 job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect refreshed market data,
# matching a scheduled market-data refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H33").Value = 839.4737
$ws.Range("I33").Value = 647.4
$ws.Range("K33").Value = 647.4
$ws.Range("M33").Value = -418.4

$ws.Range("H80").Value = 43065.12
$ws.Range("I80").Value = 70499.83
$ws.Range("J80").Value = 3284.8
$ws.Range("K80").Value = 211499.49
$ws.Range("L80").Value = 9854.400000000001
$ws.Range("M80").Value = -210501.49
$ws.Range("N80").Value = -11850.4

$ws.Range("H83").Value = 43065.12
$ws.Range("I83").Value = 70499.83
$ws.Range("J83").Value = 3284.8
$ws.Range("K83").Value = 634498.47
$ws.Range("L83").Value = 29563.2
$ws.Range("M83").Value = -629506.47
$ws.Range("N83").Value = -39547.2

$ws.Range("H106").Value = 5149333.5
$ws.Range("I106").Value = 7265971.5
$ws.Range("K106").Value = 7265971.5
$ws.Range("M106").Value = -7265340.5

$ws.Range("H107").Value = 6477.68
$ws.Range("I107").Value = 6545
$ws.Range("J107").Value = 6124.25
$ws.Range("K107").Value = 6545
$ws.Range("L107").Value = 6124.25
$ws.Range("M107").Value = -4625
$ws.Range("N107").Value = -9964.25

$ws.Range("H112").Value = 68391.2
$ws.Range("J112").Value = 85264.086
$ws.Range("L112").Value = 255792.258
$ws.Range("N112").Value = -258008.258

$ws.Range("H118").Value = 720.9
$ws.Range("I118").Value = 634.8889
$ws.Range("J118").Value = 1495
$ws.Range("K118").Value = 1904.6667
$ws.Range("L118").Value = 4485
$ws.Range("M118").Value = -247.6667000000002
$ws.Range("N118").Value = -7799

$ws.Range("H132").Value = 3720.4878
$ws.Range("I132").Value = 3905.9553
$ws.Range("J132").Value = 2892.0667
$ws.Range("K132").Value = 11717.8659
$ws.Range("L132").Value = 8676.2001
$ws.Range("M132").Value = -9187.865900000001
$ws.Range("N132").Value = -13736.2001

$ws.Range("H137").Value = 6734.2793
$ws.Range("I137").Value = 7946.636
$ws.Range("K137").Value = 23839.908
$ws.Range("M137").Value = -21289.908

$ws.Range("H138").Value = 2544.0435
$ws.Range("I138").Value = 1244.2258
$ws.Range("J138").Value = 5230.3335
$ws.Range("K138").Value = 3732.6774
$ws.Range("L138").Value = 15691.0005
$ws.Range("M138").Value = 1407.3226
$ws.Range("N138").Value = -25971.0005

$ws = $wb.Worksheets("ARM")
$ws.Range("H32").Value = 5940.486
$ws.Range("I32").Value = 6305.5156
$ws.Range("K32").Value = 6305.5156
$ws.Range("M32").Value = -6018.5156

$ws.Range("H45").Value = 11583.167
$ws.Range("I45").Value = 30000
$ws.Range("K45").Value = 30000
$ws.Range("M45").Value = -29623

$ws.Range("H61").Value = 4165.154
$ws.Range("I61").Value = 4282.4565
$ws.Range("K61").Value = 4282.4565
$ws.Range("M61").Value = -4070.4565

$ws.Range("H74").Value = 6201.077
$ws.Range("I74").Value = 1951.875
$ws.Range("K74").Value = 1951.875
$ws.Range("M74").Value = -1077.875

$ws.Range("H77").Value = 6201.077
$ws.Range("I77").Value = 1951.875
$ws.Range("K77").Value = 9759.375
$ws.Range("M77").Value = -5391.375

$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws.Range("H136").Value = 4165.154
$ws.Range("I136").Value = 4282.4565
$ws.Range("K136").Value = 12847.3695
$ws.Range("M136").Value = -10297.3695

$ws.Range("H139").Value = 102396.71
$ws.Range("J139").Value = 102396.71
$ws.Range("L139").Value = 102396.71
$ws.Range("N139").Value = -112676.71

$ws = $wb.Worksheets("BSM")
$ws.Range("H94").Value = 28223.637
$ws.Range("I94").Value = 30612.223
$ws.Range("K94").Value = 30612.223
$ws.Range("M94").Value = -30161.223

$ws.Range("H105").Value = 4227.294
$ws.Range("I105").Value = 2217.7
$ws.Range("K105").Value = 2217.7
$ws.Range("M105").Value = -470.6999999999998

$ws = $wb.Worksheets("CRP")
$ws.Range("H31").Value = 4690.4443
$ws.Range("I31").Value = 3899
$ws.Range("K31").Value = 3899
$ws.Range("M31").Value = -3604

$ws.Range("H34").Value = 4690.4443
$ws.Range("I34").Value = 3899
$ws.Range("K34").Value = 3899
$ws.Range("M34").Value = -3697

$ws.Range("H58").Value = 2032.2
$ws.Range("I58").Value = 885.94116
$ws.Range("J58").Value = 3531.1538
$ws.Range("K58").Value = 885.94116
$ws.Range("L58").Value = 3531.1538
$ws.Range("M58").Value = -682.94116
$ws.Range("N58").Value = -3937.1538

$ws.Range("H86").Value = 16000
$ws.Range("I86").Value = 10000
$ws.Range("K86").Value = 10000
$ws.Range("M86").Value = -8877

$ws.Range("H89").Value = 16000
$ws.Range("I89").Value = 10000
$ws.Range("K89").Value = 50000
$ws.Range("M89").Value = -44384

$ws.Range("H132").Value = 23042.895
$ws.Range("I132").Value = 1948.1538
$ws.Range("K132").Value = 5844.4614
$ws.Range("M132").Value = -3314.4614

$ws.Range("H134").Value = 1822.2424
$ws.Range("I134").Value = 1526.7018
$ws.Range("J134").Value = 3694
$ws.Range("K134").Value = 4580.1054
$ws.Range("L134").Value = 11082
$ws.Range("M134").Value = -2045.1054
$ws.Range("N134").Value = -16152

$ws.Range("H136").Value = 2032.2
$ws.Range("I136").Value = 885.94116
$ws.Range("J136").Value = 3531.1538
$ws.Range("K136").Value = 2657.82348
$ws.Range("L136").Value = 10593.4614
$ws.Range("M136").Value = -107.82348
$ws.Range("N136").Value = -15693.4614

$ws.Range("H138").Value = 169988.5
$ws.Range("J138").Value = 169988.5
$ws.Range("L138").Value = 169988.5
$ws.Range("N138").Value = -180268.5

$ws = $wb.Worksheets("CUL")
$ws.Range("H7").Value = 32.333332
$ws.Range("I7").Value = 47.5
$ws.Range("J7").Value = 2
$ws.Range("K7").Value = 142.5
$ws.Range("L7").Value = 6
$ws.Range("M7").Value = -30.5
$ws.Range("N7").Value = -230

$ws.Range("H61").Value = 2132.75
$ws.Range("I61").Value = 2644.1667
$ws.Range("J61").Value = 598.5
$ws.Range("K61").Value = 7932.500100000001
$ws.Range("L61").Value = 1795.5
$ws.Range("M61").Value = -7717.500100000001
$ws.Range("N61").Value = -2225.5

$ws.Range("H68").Value = 8747.706
$ws.Range("I68").Value = 1625
$ws.Range("J68").Value = 12632.818
$ws.Range("K68").Value = 4875
$ws.Range("L68").Value = 37898.454
$ws.Range("M68").Value = -4064
$ws.Range("N68").Value = -39520.454

$ws.Range("H71").Value = 8747.706
$ws.Range("I71").Value = 1625
$ws.Range("J71").Value = 12632.818
$ws.Range("K71").Value = 14625
$ws.Range("L71").Value = 113695.362
$ws.Range("M71").Value = -10569
$ws.Range("N71").Value = -121807.362

$ws.Range("H92").Value = 226.77777
$ws.Range("J92").Value = 206.25
$ws.Range("L92").Value = 618.75
$ws.Range("N92").Value = -3114.75

$ws.Range("H107").Value = 702.275
$ws.Range("J107").Value = 780.7059
$ws.Range("L107").Value = 2342.1177
$ws.Range("N107").Value = -6182.117700000001

$ws.Range("H113").Value = 3692.45
$ws.Range("I113").Value = 776.55554
$ws.Range("J113").Value = 4539
$ws.Range("K113").Value = 2329.66662
$ws.Range("L113").Value = 13617
$ws.Range("M113").Value = -159.66662
$ws.Range("N113").Value = -17957

$ws = $wb.Worksheets("GSM")
$ws.Range("H2").Value = 796.0454999999999
$ws.Range("I2").Value = 990.2941
$ws.Range("K2").Value = 990.2941
$ws.Range("M2").Value = -877.2941

$ws.Range("H11").Value = 11388797
$ws.Range("I11").Value = 12240388
$ws.Range("K11").Value = 12240388
$ws.Range("M11").Value = -12240249

$ws.Range("H18").Value = 23500000
$ws.Range("J18").Value = 20000000
$ws.Range("L18").Value = 20000000
$ws.Range("N18").Value = -20000586

$ws.Range("H43").Value = 17598.7
$ws.Range("I43").Value = 17598.7
$ws.Range("K43").Value = 17598.7
$ws.Range("M43").Value = -17447.7

$ws.Range("H102").Value = 7275.1714
$ws.Range("I102").Value = 7763.3335
$ws.Range("K102").Value = 7763.3335
$ws.Range("M102").Value = -6141.3335

$ws = $wb.Worksheets("LTW")
$ws.Range("H7").Value = 26632.809
$ws.Range("J7").Value = 7970.4116
$ws.Range("L7").Value = 7970.4116
$ws.Range("N7").Value = -8194.411599999999

$ws.Range("H40").Value = 17283.334
$ws.Range("I40").Value = 19914.87
$ws.Range("J40").Value = 12627.538
$ws.Range("K40").Value = 19914.87
$ws.Range("L40").Value = 12627.538
$ws.Range("M40").Value = -19778.87
$ws.Range("N40").Value = -12899.538

$ws.Range("H122").Value = 11099.818
$ws.Range("I122").Value = 12637.25
$ws.Range("K122").Value = 37911.75
$ws.Range("M122").Value = -35461.75

$ws.Range("H126").Value = 26632.809
$ws.Range("J126").Value = 7970.4116
$ws.Range("L126").Value = 23911.2348
$ws.Range("N126").Value = -28851.2348

$ws.Range("H132").Value = 251184.55
$ws.Range("I132").Value = 333295.25
$ws.Range("J132").Value = 4852.467
$ws.Range("K132").Value = 999885.75
$ws.Range("L132").Value = 14557.401
$ws.Range("M132").Value = -997355.75
$ws.Range("N132").Value = -19617.401

$ws.Range("H136").Value = 3615.2693
$ws.Range("I136").Value = 1465.4
$ws.Range("J136").Value = 10781.5
$ws.Range("K136").Value = 4396.200000000001
$ws.Range("L136").Value = 32344.5
$ws.Range("M136").Value = -1846.200000000001
$ws.Range("N136").Value = -37444.5

$ws = $wb.Worksheets("WVR")
$ws.Range("H132").Value = 10845.606
$ws.Range("I132").Value = 12371.083
$ws.Range("J132").Value = 6777.6665
$ws.Range("K132").Value = 37113.249
$ws.Range("L132").Value = 20332.9995
$ws.Range("M132").Value = -34583.249
$ws.Range("N132").Value = -25392.9995

$ws.Range("H136").Value = 221264.77
$ws.Range("I136").Value = 270728.47
$ws.Range("K136").Value = 812185.4099999999
$ws.Range("M136").Value = -809635.4099999999
